# KIBON-1758 add angebotGemeinden to excel template and converter
#
# Inserts a new "angebotGemeinden" column between the existing
# "gemeinden" column (G) and the "institution" column (old H, now I)
# on the Data sheet, with a title cell in row 5 and a placeholder
# cell in row 6, shifting every following column one position to the
# right (old H:Q -> new I:R).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Insert a new column at H; this shifts columns H:Q to I:R and carries
# the row-5/row-6 cell styles (header fill+border / data border) along
# with it, exactly like dragging a column-insert in the Excel UI.
$ws.Columns("H").Insert()

# New header (row 5) + placeholder (row 6) cells for the inserted column.
$ws.Range("H5").Value2 = "{angebotGemeindenTitle}"
$ws.Range("H6").Value2 = "{angebotGemeinden}"

# Give the new column its own width instead of inheriting column G's.
# (ColumnWidth is quantized to the sheet's pixel grid, same as in real
# Excel, so 22 lands on the closest achievable stored width to the
# template's 22.85546875.)
$ws.Columns("H").ColumnWidth = 22
